# "Opciones de formato para el archivo de salida"
# The question/answer table (A2:F9, 8 rows under the header in A1:F1) gets
# copied three more times below itself, stacking the same block at rows
# 10-17, 18-25 and 26-33. This mirrors what happens when a user selects the
# data block and pastes it repeatedly underneath to extend the question bank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$source = $ws.Range("A2:F9")

$source.Copy()
$ws.Range("A10").PasteSpecial()

$source.Copy()
$ws.Range("A18").PasteSpecial()

$source.Copy()
$ws.Range("A26").PasteSpecial()

# Leave the selection where the user would land next: right below the
# last pasted block (row 34 is blank, so the active cell moves to A35).
$ws.Range("A35").Select()
